# DDF2 - add soft Assertions with screenshot
# Adds a new "OpenAccountTest" worksheet (test data for the OpenAccount
# test case) right after the existing "AddCustomerTest" sheet, with a
# small customer/currency table, and leaves that new sheet as the active
# tab/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, positioned right after AddCustomerTest.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

# Header row + data row for the new OpenAccountTest data sheet.
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Alex"
$ws2.Range("B2").Value = "dollars"

# Reuse the look of AddCustomerTest's header/data rows for the new sheet.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2:B2").Copy()
$ws2.Range("A2:B2").PasteSpecial(-4122)

# AddCustomerTest's selection resets to A1 and is no longer the active tab;
# OpenAccountTest becomes the active tab with B3 selected.
$ws1.Range("A1").Select() | Out-Null
$ws2.Range("B3").Select() | Out-Null
